$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: "Time in session" - fix punctuation in the time values
$ws.Range("B7").Value = "988 hrs., 30'"
$ws.Range("C7").Value = "890 hrs., 25'"

# Row 8: "Pages of proceedings" - change from text (S10981/H10976) to numeric values
$ws.Range("B8").Value = 10981
$ws.Range("C8").Value = 10976

# Row 9: "Extension of remarks" -> "Extensions of remarks" (pluralized), and
# change C9 from text (E2386) to numeric value
$ws.Range("A9").Value = "Extensions of remarks"
$ws.Range("C9").Value = 2386
